# daily auto push: 2026-02-13 19:06 UTC
#
# A new observation row for 2026/02/14 (Sat) slipped into the log between
# the 2026/02/13 entries and the 2026/12/29 entries. Insert it at row 806
# and let every following row shift down by one (dimension grows from
# D847 to D848).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 806..847 down to 807..848, leaving row 806 blank for the new entry.
$ws.Rows.Item(806).Insert()

# Column A holds a date formatted as plain text (e.g. "2025/01/01"), not a
# real date value. Using a leading apostrophe forces Excel to store the
# literal text instead of auto-converting it to a date serial number, and
# ClearFormats() drops the "entered as text" quote-prefix formatting so the
# cell ends up with the same default (unstyled) look as its neighbours.
$ws.Cells.Item(806, 1).Value = "'2026/02/14"
$ws.Cells.Item(806, 1).ClearFormats()

$ws.Cells.Item(806, 2).Value = "土"
$ws.Cells.Item(806, 3).Value = 1
$ws.Cells.Item(806, 4).Value = 201
